# Update the "Metadata" sheet of the ValueSet-ConsentExchangePolicies workbook:
#  - refresh the Date value
#  - replace the placeholder Contact value with the real contact info
#  - insert a new "Jurisdiction" property row (with an empty value) right
#    after "Contact", which pushes Description/Purpose/Copyright/Immutable
#    down by one row (Copyright's value stays empty, i.e. removed)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Date property (row 8) gets a new timestamp.
$ws.Range("B8").Value = "2024-09-09T14:48:24-05:00"

# 2. Contact property (row 10) gets real contact info instead of the
#    "No display for ContactDetail" placeholder.
$ws.Range("B10").Value = "MITRE, Inc (https://github.com/awatson1978/us-state-profiles)"

# 3. Insert a new row right after "Contact" (row 10) for "Jurisdiction".
$ws.Rows.Item(11).Insert()

# Copy the formatting of the row above onto the newly-inserted row so it
# keeps the same borders/alignment as the rest of the property table.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's values: Property = "Jurisdiction", Value = blank.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
